$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, pushing existing rows 30..132 down to 31..133.
# Excel's default Insert() behavior copies formatting from the row above,
# which correctly preserves the date format (style) for column D.
$ws.Rows("30:30").Insert()

# Populate the newly inserted row 30 with the new record's data.
$ws.Cells.Item(30, 1).Value = 1
$ws.Cells.Item(30, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(30, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(30, 4).Value = "2022-12-15"
$ws.Cells.Item(30, 5).Value = 15
$ws.Cells.Item(30, 6).Value = "Fruta"
$ws.Cells.Item(30, 7).Value = 100102
$ws.Cells.Item(30, 8).Value = "Cítricos"
$ws.Cells.Item(30, 9).Value = 100102004
$ws.Cells.Item(30, 10).Value = "Mandarina"
$ws.Cells.Item(30, 11).Value = "Murcott"
$ws.Cells.Item(30, 12).Value = "Tercera"
$ws.Cells.Item(30, 13).Value = 350
$ws.Cells.Item(30, 14).Value = 14000
$ws.Cells.Item(30, 15).Value = 15000
$ws.Cells.Item(30, 16).Value = 14429
$ws.Cells.Item(30, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(30, 18).Value = "Región Metropolitana"
$ws.Cells.Item(30, 19).Value = 721
$ws.Cells.Item(30, 20).Value = 20
